$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (Control 10)
$ws.Range("D2").Value = 0.00002361557811981688
$ws.Range("E2").Value = 0.00002361557811981688

# Row 3 (Control 5)
$ws.Range("D3").Value = 0.9993276578237307
$ws.Range("E3").Value = 0.9993276578237307

# Row 4 (MDD 42)
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = 0.002848182626287424
$ws.Range("E4").Value = 0.9971518173737126

# Row 6 (MDD 20)
$ws.Range("D6").Value = 0.9982600664929326
$ws.Range("E6").Value = 0.001739933507067359

# Row 7 (MDD 51)
$ws.Range("D7").Value = 0.00001989479141991188
$ws.Range("E7").Value = 0.9999801052085801

# Row 8 (MDD 40)
$ws.Range("D8").Value = 0.0000000000000230213536446217
$ws.Range("E8").Value = 0.999999999999977
$ws.Range("F8").Value = 7.913569927215576
$ws.Range("G8").Value = 0.4285714285714285
